$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value2 = "KAOGExp"
$ws.Range("B2").Value2 = 0.99
$ws.Range("C2").Value2 = 5.161616161616162
$ws.Range("D2").Value2 = 2.199146899475081
$ws.Range("E2").Value2 = 2.163708741019674
$ws.Range("F2").Value2 = 1

$ws.Range("A3").Value2 = "cruds"
$ws.Range("B3").Value2 = 1
$ws.Range("C3").Value2 = 5.03
$ws.Range("D3").Value2 = 2.056648763949894
$ws.Range("E3").Value2 = 1.516368121760769
$ws.Range("F3").Value2 = 0.9128466555946751

$ws.Range("A4").Value2 = "wachter"
$ws.Range("B4").Value2 = 0.77
$ws.Range("C4").Value2 = 3.56
$ws.Range("D4").Value2 = 5.208734676873722
$ws.Range("E4").Value2 = 100.1133167941762
$ws.Range("F4").Value2 = 1.70617793086653

$ws.Range("A5").Value2 = "face-knn"
$ws.Range("B5").Value2 = 1
$ws.Range("C5").Value2 = 5.12
$ws.Range("D5").Value2 = 2.806598448043185
$ws.Range("E5").Value2 = 2.476218512671286
$ws.Range("F5").Value2 = 0.9553778677462887

$ws.Range("A6").Value2 = "revise"
$ws.Range("B6").Value2 = 1
$ws.Range("C6").Value2 = 4.83
$ws.Range("D6").Value2 = 1.824670988148446
$ws.Range("E6").Value2 = 1.224971012286851
$ws.Range("F6").Value2 = 0.8813361167147559

$ws.Range("A7").Value2 = "cem"
$ws.Range("B7").Value2 = 1
$ws.Range("C7").Value2 = 4
$ws.Range("D7").Value2 = 1.389180626876773
$ws.Range("E7").Value2 = 1.0909338597036
$ws.Range("F7").Value2 = 0.8352631578947368

$ws.Range("A8").Value2 = "cem-vae"
$ws.Range("B8").Value2 = 1
$ws.Range("C8").Value2 = 3.98
$ws.Range("D8").Value2 = 1.357806435920778
$ws.Range("E8").Value2 = 1.052866545500113
$ws.Range("F8").Value2 = 0.8352631578947368

$ws.Range("A9").Value2 = "dice"
$ws.Range("B9").Value2 = 1
$ws.Range("C9").Value2 = 2.44
$ws.Range("D9").Value2 = 1.761591970310391
$ws.Range("E9").Value2 = 1.533711892708396
$ws.Range("F9").Value2 = 0.9297311403508772

$ws.Range("A10").Value2 = "face-epsilon"
$ws.Range("B10").Value2 = 1
$ws.Range("C10").Value2 = 4.97
$ws.Range("D10").Value2 = 2.725550185560054
$ws.Range("E10").Value2 = 2.405854398579997
$ws.Range("F10").Value2 = 0.9589381916329285

$ws.Range("A11").Value2 = "clue"
$ws.Range("B11").Value2 = 1
$ws.Range("C11").Value2 = 4.83
$ws.Range("D11").Value2 = 2.076470991312613
$ws.Range("E11").Value2 = 1.454493131918217
$ws.Range("F11").Value2 = 0.8739165748862281

$ws.Range("A12").Value2 = "ar"
$ws.Range("B12").Value2 = 0.53
$ws.Range("C12").Value2 = 0.9
$ws.Range("D12").Value2 = 0.8163907499999999
$ws.Range("E12").Value2 = 0.8045457230575311
$ws.Range("F12").Value2 = 0.5073085394736842

$ws.Range("A13").Value2 = "cchvae"
$ws.Range("B13").Value2 = 1
$ws.Range("C13").Value2 = 5.73
$ws.Range("D13").Value2 = 2.803325219740157
$ws.Range("E13").Value2 = 2.139333228311118
$ws.Range("F13").Value2 = 0.9824790439793938

$ws.Range("A14").Value2 = "gs"
$ws.Range("B14").Value2 = 1
$ws.Range("C14").Value2 = 3.76
$ws.Range("D14").Value2 = 0.9528221717893516
$ws.Range("E14").Value2 = 0.790185794911023
$ws.Range("F14").Value2 = 0.7527914303274704
